$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 115, shifting the existing rows 115-127 down to 116-128.
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row 115 with the new weekly price record.
$ws.Cells.Item(115, 1).Value = 7
$ws.Cells.Item(115, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(115, 3).Value = "Ñuble"
$ws.Cells.Item(115, 4).Value = 44474
$ws.Cells.Item(115, 5).Value = 16
$ws.Cells.Item(115, 6).Value = 100112003
$ws.Cells.Item(115, 7).Value = "Ajo"
$ws.Cells.Item(115, 8).Value = "Chino"
$ws.Cells.Item(115, 9).Value = "Primera"
$ws.Cells.Item(115, 10).Value = 60
$ws.Cells.Item(115, 11).Value = 16000
$ws.Cells.Item(115, 12).Value = 17000
$ws.Cells.Item(115, 13).Value = 16500
$ws.Cells.Item(115, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(115, 15).Value = "China"
$ws.Cells.Item(115, 16).Value = 1650
$ws.Cells.Item(115, 17).Value = 10
$ws.Cells.Item(115, 18).Value = "Hortaliza"
